# sprint_backlog.xlsx - complete M3 sprint backlog sheet
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Spring 1 (M1)
$ws2 = $wb.Worksheets.Item(2)   # Sprint 2 (M2)
$ws3 = $wb.Worksheets.Item(3)   # Sprint 3 (M3)

# ---------------------------------------------------------------------------
# Sprint 2 (M2): a few task rows got taller (wrapped description text grew)
# ---------------------------------------------------------------------------
2..7 | ForEach-Object { $ws2.Rows.Item($_).RowHeight = 31.5 }
$ws2.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sprint 3 (M3): build out the backlog table that was previously empty
# ---------------------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 32.25

function Fill-Row($r, $a, $b, $c, $d) {
    $ws3.Range("A$r").Value = $a
    $ws3.Range("B$r").Value = $b
    $ws3.Range("C$r").Value = $c
    $ws3.Range("D$r").Value = $d
}

# Header row
Fill-Row 1 "Tasks" "Responsible" "Status" 1
$ws3.Range("E1").Value = 2
$ws3.Range("F1").Value = 3
$ws3.Range("G1").Value = 4

# Data rows -- entered in this order so new shared strings land at the same
# indices as in the authored workbook (use-case tasks first, then the
# remaining M3 tasks)
Fill-Row 2 "Create Context Diagram" "Pranil" "Not Started" 1
Fill-Row 8 "Collate the best use cases" "Hunter" "Not Started" 0.5
Fill-Row 3 "Brainstorm 10 use cases" "Hunter" "Not Started" 0.5
Fill-Row 4 "Brainstorm 10 use cases" "Bhavesh" "Not Started" 0.5
Fill-Row 5 "Brainstorm 10 use cases" "Stephen" "Not Started" 0.5
Fill-Row 6 "Brainstorm 10 use cases" "Pranil" "Not Started" 0.5
Fill-Row 7 "Brainstorm 10 use cases" "Naman" "Not Started" 0.5
Fill-Row 13 "Handle player object creation and game configuration dialog mechanics (5 & 6)" "Naman" "Not Started" 1
Fill-Row 9 "Create Space Trader FX project (1)" "Hunter" "Not Started" 1
Fill-Row 10 "Design Welcome Screen UI (2)" "Stephen" "Not Started" 1
Fill-Row 11 "Implement pressing ""New Game"" to access the game configuration dialog (3)" "Bhavesh" "Not Started" 1
Fill-Row 12 "Implement character creation (name + skill points) (4)" "Pranil" "Not Started" 1

# Header formatting: bold + wrap text (matches the style used on the other
# sheets' header rows)
$ws3.Range("A1:G1").WrapText = $true
$ws3.Range("A1:G1").Font.Bold = $true

# Data cell formatting: wrap text on the descriptive columns only
$ws3.Range("A2:C13").WrapText = $true

# Rows whose wrapped text needed extra height
$ws3.Rows.Item(11).RowHeight = 47.25
$ws3.Rows.Item(12).RowHeight = 31.5
$ws3.Rows.Item(13).RowHeight = 47.25

$ws3.Range("F12").Select() | Out-Null

# ---------------------------------------------------------------------------
# Make Sprint 3 (M3) the active tab
# ---------------------------------------------------------------------------
$ws3.Activate() | Out-Null

Write-Host "Sprint backlog M3 populated"
